$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 267.53333
$ws.Range("I33").Value = 173.88889
$ws.Range("J33").Value = 408
$ws.Range("K33").Value = 173.88889
$ws.Range("L33").Value = 408
$ws.Range("M33").Value = 55.11111
$ws.Range("N33").Value = -866

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 146.23077
$ws.Range("I55").Value = 125.6
$ws.Range("J55").Value = 215
$ws.Range("K55").Value = 125.6
$ws.Range("L55").Value = 215
$ws.Range("M55").Value = 88.4
$ws.Range("N55").Value = -643

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 37749
$ws.Range("J105").Value = 37749
$ws.Range("L105").Value = 37749
$ws.Range("N105").Value = -44737

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2670.1428
$ws.Range("I137").Value = 2372.75
$ws.Range("J137").Value = 3066.6667
$ws.Range("K137").Value = 7118.25
$ws.Range("L137").Value = 9200.000100000001
$ws.Range("M137").Value = -4568.25
$ws.Range("N137").Value = -14300.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 900
$ws.Range("I74").Value = 900
$ws.Range("K74").Value = 900
$ws.Range("M74").Value = -26

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 900
$ws.Range("I77").Value = 900
$ws.Range("K77").Value = 4500
$ws.Range("M77").Value = -132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 26349.8
$ws.Range("J106").Value = 26349.8
$ws.Range("L106").Value = 26349.8
$ws.Range("N106").Value = -28873.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1622.2
$ws.Range("I132").Value = 1527.75
$ws.Range("K132").Value = 4583.25
$ws.Range("M132").Value = -2053.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1117.0714
$ws.Range("I64").Value = 1571.75
$ws.Range("J64").Value = 935.2
$ws.Range("K64").Value = 1571.75
$ws.Range("L64").Value = 935.2
$ws.Range("M64").Value = -1346.75
$ws.Range("N64").Value = -1385.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1117.0714
$ws.Range("I67").Value = 1571.75
$ws.Range("J67").Value = 935.2
$ws.Range("K67").Value = 1571.75
$ws.Range("L67").Value = 935.2
$ws.Range("M67").Value = -791.75
$ws.Range("N67").Value = -2495.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2178.4546
$ws.Range("I94").Value = 2440.3333
$ws.Range("K94").Value = 2440.3333
$ws.Range("M94").Value = -1989.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1241.8572
$ws.Range("I107").Value = 1019.8
$ws.Range("K107").Value = 1019.8
$ws.Range("M107").Value = 900.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6666.1665
$ws.Range("I31").Value = 4999
$ws.Range("J31").Value = 8333.333
$ws.Range("K31").Value = 4999
$ws.Range("L31").Value = 8333.333
$ws.Range("M31").Value = -4704
$ws.Range("N31").Value = -8923.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6666.1665
$ws.Range("I34").Value = 4999
$ws.Range("J34").Value = 8333.333
$ws.Range("K34").Value = 4999
$ws.Range("L34").Value = 8333.333
$ws.Range("M34").Value = -4797
$ws.Range("N34").Value = -8737.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 13500
$ws.Range("J43").Value = 13500
$ws.Range("L43").Value = 13500
$ws.Range("N43").Value = -13868

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 13500
$ws.Range("J101").Value = 13500
$ws.Range("L101").Value = 13500
$ws.Range("N101").Value = -19990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1799.2
$ws.Range("I107").Value = 1999.3334
$ws.Range("J107").Value = 1499
$ws.Range("K107").Value = 1999.3334
$ws.Range("L107").Value = 1499
$ws.Range("M107").Value = -79.33339999999998
$ws.Range("N107").Value = -5339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3870
$ws.Range("I134").Value = 3806.4443
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 11419.3329
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -8884.332900000001
$ws.Range("N134").Value = -20112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 663.375
$ws.Range("J23").Value = 584.5
$ws.Range("L23").Value = 1753.5
$ws.Range("N23").Value = -2223.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 96348.75
$ws.Range("J37").Value = 96348.75
$ws.Range("L37").Value = 289046.25
$ws.Range("N37").Value = -289270.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 464.16666
$ws.Range("I97").Value = 253.75
$ws.Range("J97").Value = 885
$ws.Range("K97").Value = 761.25
$ws.Range("L97").Value = 2655
$ws.Range("M97").Value = -265.25
$ws.Range("N97").Value = -3647

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5001980
$ws.Range("I122").Value = 5001980
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15005940
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15003490
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4199.75
$ws.Range("I132").Value = 3899.5
$ws.Range("K132").Value = 11698.5
$ws.Range("M132").Value = -9168.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1875.4375
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3333.3333
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 3333.3333
$ws.Range("M2").Value = 111
$ws.Range("N2").Value = -3557.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 24398.6
$ws.Range("J103").Value = 24398.6
$ws.Range("L103").Value = 24398.6
$ws.Range("N103").Value = -26742.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 41500
$ws.Range("J104").Value = 41500
$ws.Range("L104").Value = 41500
$ws.Range("N104").Value = -48488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 29428.572
$ws.Range("J129").Value = 29428.572
$ws.Range("L129").Value = 29428.572
$ws.Range("N129").Value = -39428.572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2442.1428
$ws.Range("I132").Value = 2515.8333
$ws.Range("K132").Value = 7547.499899999999
$ws.Range("M132").Value = -5017.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1552.1818
$ws.Range("I136").Value = 1008.2222
$ws.Range("K136").Value = 3024.6666
$ws.Range("M136").Value = -474.6666
